# Add Employee to db feature rdy
# Opmann_Time_recording_log.xlsx - "Nädal 4" sheet (4th worksheet)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(4)

# Row 10: mark activity H10 as left-aligned (new style)
$ws.Range("H10").HorizontalAlignment = -4131

# Row 11 (entry #5): fill in date/time worked, activity "Prep." and comment "ASP.NET MVC"
$ws.Range("B11").Value = 43520
$ws.Range("C11").Value = 0.29166666666666669
$ws.Range("D11").Value = 0.36805555555555558
$ws.Range("G11").Value = "Prep."
$ws.Range("H11").Value = "ASP.NET MVC"
$ws.Range("H11").HorizontalAlignment = -4131

# Row 12 (entry #6): fill in time worked and comment "ASP.NET MVC"
$ws.Range("C12").Value = 0.4375
$ws.Range("D12").Value = 0.5
$ws.Range("H12").Value = "ASP.NET MVC"
$ws.Range("H12").HorizontalAlignment = -4131

# Update the active selection left on the sheet
[void]$ws.Range("D13").Select()
